$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename move ids so the numeric suffix is zero-padded to 2 digits
# (z0bug.move_1 -> z0bug.move_01, z0bug.move_1_1 -> z0bug.move_01_1, etc.)
# Column A holds the line id (move_X_Y), column B holds the move id (move_X)
$updates = @(
    @{ Row = 2;  A = "z0bug.move_01_1"; B = "z0bug.move_01" },
    @{ Row = 3;  A = "z0bug.move_01_2"; B = "z0bug.move_01" },
    @{ Row = 4;  A = "z0bug.move_02_1"; B = "z0bug.move_02" },
    @{ Row = 5;  A = "z0bug.move_02_2"; B = "z0bug.move_02" },
    @{ Row = 6;  A = "z0bug.move_03_1"; B = "z0bug.move_03" },
    @{ Row = 7;  A = "z0bug.move_03_2"; B = "z0bug.move_03" },
    @{ Row = 8;  A = "z0bug.move_03_3"; B = "z0bug.move_03" },
    @{ Row = 9;  A = "z0bug.move_04_1"; B = "z0bug.move_04" },
    @{ Row = 10; A = "z0bug.move_04_2"; B = "z0bug.move_04" },
    @{ Row = 11; A = "z0bug.move_05_1"; B = "z0bug.move_05" },
    @{ Row = 12; A = "z0bug.move_05_2"; B = "z0bug.move_05" },
    @{ Row = 13; A = "z0bug.move_06_1"; B = "z0bug.move_06" },
    @{ Row = 14; A = "z0bug.move_06_2"; B = "z0bug.move_06" },
    @{ Row = 15; A = "z0bug.move_07_1"; B = "z0bug.move_07" },
    @{ Row = 16; A = "z0bug.move_07_2"; B = "z0bug.move_07" },
    @{ Row = 17; A = "z0bug.move_07_3"; B = "z0bug.move_07" },
    @{ Row = 18; A = "z0bug.move_08_1"; B = "z0bug.move_08" },
    @{ Row = 19; A = "z0bug.move_08_2"; B = "z0bug.move_08" },
    @{ Row = 20; A = "z0bug.move_09_1"; B = "z0bug.move_09" },
    @{ Row = 21; A = "z0bug.move_09_2"; B = "z0bug.move_09" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 2).Value = $u.B
}

# Update the active selection on the sheet (was E11, now C5)
$ws.Range("C5").Select()

# Widen columns A and B slightly (15.68 -> 17.51, 14.03 -> 14.59 in the saved file).
# The ColumnWidth property is quantized to pixel units by the engine, so the values
# below are chosen to round-trip to the closest representable width.
$ws.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334
